$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2860.6155
$ws.Range("J17").Value = 2860.6155
$ws.Range("L17").Value = 8581.8465
$ws.Range("N17").Value = -8917.8465
$ws.Range("H28").Value = 2300
$ws.Range("I28").Value = 1900
$ws.Range("J28").Value = 3500
$ws.Range("K28").Value = 1900
$ws.Range("L28").Value = 3500
$ws.Range("M28").Value = -1415
$ws.Range("N28").Value = -4470
$ws.Range("H41").Value = 694
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H86").Value = 6000
$ws.Range("I86").Value = 6000
$ws.Range("K86").Value = 6000
$ws.Range("H89").Value = 6000
$ws.Range("I89").Value = 6000
$ws.Range("K89").Value = 30000
$ws.Range("H96").Value = 100000000
$ws.Range("J96").Value = 100000000
$ws.Range("L96").Value = 300000000
$ws.Range("H98").Value = 1267.3334
$ws.Range("I98").Value = 1267.3334
$ws.Range("K98").Value = 1267.3334
$ws.Range("M98").Value = 230.6666
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("H116").Value = 5816.5
$ws.Range("I116").Value = 5474.75
$ws.Range("K116").Value = 5474.75
$ws.Range("M116").Value = -2032.75
$ws.Range("H122").Value = 1267.3334
$ws.Range("I122").Value = 1267.3334
$ws.Range("K122").Value = 3802.0002
$ws.Range("M122").Value = -1352.0002
$ws.Range("H132").Value = 1457.7142
$ws.Range("I132").Value = 1457.7142
$ws.Range("K132").Value = 4373.142599999999
$ws.Range("M132").Value = -1843.142599999999
$ws.Range("M86").Value = -4877
$ws.Range("M89").Value = -24384
$ws.Range("N96").Value = -300002746
$ws.Range("N41").ClearContents()
$ws.Range("M111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 8000
$ws.Range("I11").Value = 8000
$ws.Range("K11").Value = 8000
$ws.Range("M11").Value = -7856
$ws.Range("H32").Value = 15050.052
$ws.Range("I32").Value = 14637.639
$ws.Range("K32").Value = 14637.639
$ws.Range("M32").Value = -14350.639
$ws.Range("H37").Value = 38534.5
$ws.Range("J37").Value = 38534.5
$ws.Range("L37").Value = 38534.5
$ws.Range("N37").Value = -39080.5
$ws.Range("H44").Value = 39993.332
$ws.Range("J44").Value = 39993.332
$ws.Range("L44").Value = 39993.332
$ws.Range("N44").Value = -40969.332
$ws.Range("H61").Value = 1842.963
$ws.Range("I61").Value = 1428.75
$ws.Range("K61").Value = 1428.75
$ws.Range("M61").Value = -1216.75
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H132").Value = 2703
$ws.Range("J132").Value = 5332.6665
$ws.Range("L132").Value = 15997.9995
$ws.Range("N132").Value = -21057.9995
$ws.Range("H136").Value = 1842.963
$ws.Range("I136").Value = 1428.75
$ws.Range("K136").Value = 4286.25
$ws.Range("M136").Value = -1736.25
$ws.Range("N98").ClearContents()
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 576.5
$ws.Range("J94").Value = 1337
$ws.Range("L94").Value = 1337
$ws.Range("N94").Value = -2239
$ws.Range("H105").Value = 7521.6665
$ws.Range("I105").Value = 7461.875
$ws.Range("K105").Value = 7461.875
$ws.Range("M105").Value = -5714.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1999.5
$ws.Range("I16").Value = 1999
$ws.Range("K16").Value = 1999
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H99").Value = 9407
$ws.Range("I99").Value = 9147.666999999999
$ws.Range("K99").Value = 9147.666999999999
$ws.Range("M99").Value = -7649.666999999999
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("H113").Value = 1999.5
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("H126").Value = 9407
$ws.Range("I126").Value = 9147.666999999999
$ws.Range("K126").Value = 27443.001
$ws.Range("M126").Value = -24973.001
$ws.Range("H132").Value = 4505.0435
$ws.Range("I132").Value = 4458
$ws.Range("K132").Value = 13374
$ws.Range("M132").Value = -10844
$ws.Range("M16").Value = -1712
$ws.Range("M113").Value = 171
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2800
$ws.Range("I31").Value = 2800
$ws.Range("K31").Value = 8400
$ws.Range("H56").Value = 19608.666
$ws.Range("I56").Value = 19608.666
$ws.Range("K56").Value = 19608.666
$ws.Range("M56").Value = -19078.666
$ws.Range("M31").Value = -8112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("H46").Value = 4324.2856
$ws.Range("J46").Value = 4324.2856
$ws.Range("L46").Value = 4324.2856
$ws.Range("N46").Value = -4636.2856
$ws.Range("H80").Value = 7256.5713
$ws.Range("I80").Value = 4749
$ws.Range("K80").Value = 4749
$ws.Range("M80").Value = -3751
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H83").Value = 7256.5713
$ws.Range("I83").Value = 4749
$ws.Range("K83").Value = 23745
$ws.Range("M83").Value = -18753
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7340
$ws.Range("H134").Value = 50112.5
$ws.Range("J134").Value = 50112.5
$ws.Range("L134").Value = 150337.5
$ws.Range("N134").Value = -155407.5
$ws.Range("N15").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4247.4
$ws.Range("I61").Value = 4184.25
$ws.Range("K61").Value = 4184.25
$ws.Range("M61").Value = -3982.25
$ws.Range("H113").Value = 4247.4
$ws.Range("I113").Value = 4184.25
$ws.Range("K113").Value = 4184.25
$ws.Range("M113").Value = -2014.25
$ws.Range("H136").Value = 6404.3125
$ws.Range("I136").Value = 6359.846
$ws.Range("K136").Value = 19079.538
$ws.Range("M136").Value = -16529.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 312.25
$ws.Range("I113").Value = 200
$ws.Range("J113").Value = 379.6
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 1138.8
$ws.Range("M113").Value = 1570
$ws.Range("N113").Value = -5478.8
$ws.Range("H132").Value = 2383.3914
$ws.Range("I132").Value = 1519.2941
$ws.Range("K132").Value = 4557.8823
$ws.Range("M132").Value = -2027.8823

Write-Host "Applied all updates"